$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - RandomForestClassifier
$ws.Range("C2").Value = 0.744672131147541
$ws.Range("D2").Value = 0.6926273982439175
$ws.Range("E2").Value = 0.7966029111793045
$ws.Range("F2").Value = 0.85
$ws.Range("G2").Value = 0.639344262295082
$ws.Range("H2").Value = 0.8528428093645485
$ws.Range("I2").Value = 0.85
$ws.Range("J2").Value = 0.8514190317195326

# Row 3 - XGBClassifier
$ws.Range("C3").Value = 0.7650273224043715
$ws.Range("D3").Value = 0.7132603027129654
$ws.Range("E3").Value = 0.8171273354922864
$ws.Range("F3").Value = 0.8333333333333334
$ws.Range("G3").Value = 0.6967213114754098
$ws.Range("H3").Value = 0.8710801393728222
$ws.Range("I3").Value = 0.8333333333333334
$ws.Range("J3").Value = 0.8517887563884157

# Row 4 - LogisticRegression
$ws.Range("C4").Value = 0.7674590163934426
$ws.Range("D4").Value = 0.716234660817073
$ws.Range("E4").Value = 0.8186229259398896
$ws.Range("F4").Value = 0.83
$ws.Range("G4").Value = 0.7049180327868853
$ws.Range("H4").Value = 0.8736842105263158
$ws.Range("I4").Value = 0.83
$ws.Range("J4").Value = 0.8512820512820513
